$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header fields
$ws.Range("C2").Value = "Hartmut"
# Leading apostrophe forces text storage (card number must stay a text value,
# not be auto-converted to a number) - standard Excel quote-prefix entry.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance
$ws.Range("D5").Value = "KONTOSTAND AM 16.05.2025"

# Row 6
$ws.Range("B6").Value = "18.05."
$ws.Range("C6").Value = "19.05."
$ws.Range("D6").Value = "BURGER KING Lemgo"
$ws.Range("E6").Value = "24,24-"

# Row 7
$ws.Range("B7").Value = "21.05."
$ws.Range("C7").Value = "22.05."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 48633916"
$ws.Range("E7").Value = "87,67-"

# Row 8
$ws.Range("B8").Value = "23.05."
$ws.Range("C8").Value = "24.05."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "24,97-"

# Rows 9, 10, 11 - transactions removed, clear their contents
$ws.Range("B9:E9").Value = ""
$ws.Range("B10:E10").Value = ""
$ws.Range("B11:E11").Value = ""

# Update alignment for E9 (center) and E10/E11 (right) to match the reduced-data layout
$ws.Range("E9").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E9").VerticalAlignment = -4108    # xlCenter
$ws.Range("E9").WrapText = $true

$ws.Range("E10").HorizontalAlignment = -4152 # xlRight
$ws.Range("E10").VerticalAlignment = -4108   # xlCenter
$ws.Range("E10").WrapText = $true

$ws.Range("E11").HorizontalAlignment = -4152 # xlRight
$ws.Range("E11").VerticalAlignment = -4108   # xlCenter
$ws.Range("E11").WrapText = $true

# Closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 26.05.2025"
$ws.Range("E12").Value = "136,88-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 03.06.2025"
